$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Summary  (overall account stats after Trade #52 closed)
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.8     # Current Capital
$wsSummary.Range("B4").Value = 0.9        # Total P&L $
$wsSummary.Range("B5").Value = 0.36       # Total P&L %
$wsSummary.Range("B6").Value = 50         # Total Trades
$wsSummary.Range("B8").Value = 19         # Losing Trades
$wsSummary.Range("B9").Value = 56         # Win Rate %

# ---------------------------------------------------------------
# Sheet: Strategy Status  (MarketMaking row)
# ---------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 99.8        # Capital
$wsStatus.Range("D6").Value = 21          # Trades
$wsStatus.Range("E6").Value = -0.01       # P&L $
$wsStatus.Range("F6").Value = -0.2        # P&L %
$wsStatus.Range("G6").Value = 57.14       # Win Rate %

# ---------------------------------------------------------------
# Sheet: All Trades  (Trade #52, row 53 -> closed early)
# ---------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Range("G53").Value = 0.24
$wsAllTrades.Range("H53").Value = "CLOSED"
$wsAllTrades.Range("I53").Value = -22.5806
$wsAllTrades.Range("J53").Value = -0.07000000000000001
$wsAllTrades.Range("K53").Value = 99.8
$wsAllTrades.Range("L53").Value = "early_exit"
$wsAllTrades.Range("M53").Value = 0.16

# ---------------------------------------------------------------
# Sheet: MarketMaking  (Trade #52, row 24 -> closed early)
# ---------------------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
$wsMarketMaking.Range("G24").Value = 0.24
$wsMarketMaking.Range("H24").Value = "CLOSED"
$wsMarketMaking.Range("I24").Value = -22.5806
$wsMarketMaking.Range("J24").Value = -0.07000000000000001
$wsMarketMaking.Range("K24").Value = 99.8
$wsMarketMaking.Range("P24").Value = "early_exit"
$wsMarketMaking.Range("Q24").Value = 0.16
